# FeatureMaster.xlsx - "Excel changes for feature/model"
#
# The header row is restructured:
#   - "DisplayName"          (B1) -> renamed to "FeatureName"
#   - "FeatureStatus"        (F1) -> renamed to "Status"
#   - "ModelTemplateItem"    (G1) -> removed
#   - "ItemCodeGenerationRef"(H1) -> removed
#   - "PicturePath"/"Accessory" (old I1/J1) shift left into G1/H1
#   - the stray red-font highlight on A1 ("FeatureCode") is cleared
#   - selection/view settle on the new last column (G1 = PicturePath)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename columns in place (keeps column-width metadata untouched).
$ws.Range("B1").Value = "FeatureName"
$ws.Range("F1").Value = "Status"

# Shift the trailing two columns left, onto the two columns being removed,
# then clear what used to be the last two (now-duplicate) columns.
$ws.Range("G1").Value = "PicturePath"
$ws.Range("H1").Value = "Accessory"
$ws.Range("I1:J1").ClearContents()

# Drop the direct (red) font formatting that was on A1.
$ws.Range("A1").ClearFormats()

# Reset the view: no frozen/scrolled-off left column, selection on G1.
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
[void]$ws.Range("G1").Select()
